# Insert two new weekly records right after row 295 ("Fruta / hortaliza,
# semanal") for the "Terminal Hortofrutícola Agro Chillán" - Apio sheet.
# Everything from the old row 296 onward shifts down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 296 (pushes old rows 296-315 -> 298-317).
$ws.Rows.Item(296).Resize(2).Insert()

# Shared / constant attributes for this data block.
$mercadoId   = 7
$mercado     = "Terminal Hortofrutícola Agro Chillán"
$region      = "Ñuble"
$codreg      = 16
$categoriaId = 100112017
$categoria   = "Apio"
$variedad    = "Americana (o)"
$unidad      = "`$/docena de matas"
$kgUnidades  = 6
$clasif      = "Hortaliza"

# New row 296
$r = 296
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = 45021
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $categoriaId
$ws.Cells.Item($r, 7).Value  = $categoria
$ws.Cells.Item($r, 8).Value  = $variedad
$ws.Cells.Item($r, 9).Value  = "Primera"
$ws.Cells.Item($r, 10).Value = 100
$ws.Cells.Item($r, 11).Value = 7000
$ws.Cells.Item($r, 12).Value = 7000
$ws.Cells.Item($r, 13).Value = 7000
$ws.Cells.Item($r, 14).Value = $unidad
$ws.Cells.Item($r, 15).Value = "Provincia del Elquí"
$ws.Cells.Item($r, 16).Value = 1167
$ws.Cells.Item($r, 17).Value = $kgUnidades
$ws.Cells.Item($r, 18).Value = $clasif

# New row 297
$r = 297
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = 45021
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $categoriaId
$ws.Cells.Item($r, 7).Value  = $categoria
$ws.Cells.Item($r, 8).Value  = $variedad
$ws.Cells.Item($r, 9).Value  = "Segunda"
$ws.Cells.Item($r, 10).Value = 50
$ws.Cells.Item($r, 11).Value = 5000
$ws.Cells.Item($r, 12).Value = 6000
$ws.Cells.Item($r, 13).Value = 5600
$ws.Cells.Item($r, 14).Value = $unidad
$ws.Cells.Item($r, 15).Value = "Provincia del Elquí"
$ws.Cells.Item($r, 16).Value = 933
$ws.Cells.Item($r, 17).Value = $kgUnidades
$ws.Cells.Item($r, 18).Value = $clasif
